$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.878.48'
$ws.Range("E2").Value = '  +6.40%  '

# Row 3
$ws.Range("D3").Value = '3.583.36'
$ws.Range("E3").Value = '  +5.41%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.88%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.08%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.644'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.82%  '

# Row 8
$ws.Range("D8").Value = '3.575.41'
$ws.Range("E8").Value = '  +5.57%  '

# Row 9
$ws.Range("E9").Value = '  -0.08%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.184'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.74%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.663'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.86%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000293'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.55%  '

# Row 14
$ws.Range("E14").Value = '  +5.67%  '

# Row 15
$ws.Range("D15").Value = '4.147.56'
$ws.Range("E15").Value = '  +5.34%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.74%  '

# Row 17
$ws.Range("D17").Value = '3.576.36'
$ws.Range("E17").Value = '  +4.90%  '

# Row 18
$ws.Range("D18").Value = '69.672.03'
$ws.Range("E18").Value = '  +6.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.34%  '

# Row 20
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("E21").Value = '  +4.96%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '500.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.24%  '

# Row 23
$ws.Range("E23").Value = '  +11.56%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +20.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.69%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.05%  '

# Row 27
$ws.Range("E27").Value = '  +5.02%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.10%  '

# Row 29
$ws.Range("E29").Value = '  +7.18%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.98%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '612.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.09%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '65.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.80%  '

# Row 35
$ws.Range("E35").Value = '  +7.15%  '

# Row 36
$ws.Range("D36").Value = '0.0₃0839'
$ws.Range("E36").Value = '  +12.75%  '

# Row 37
$ws.Range("E37").Value = '  +4.66%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '38.12'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.38%  '

# Row 39
$ws.Range("E39").Value = '  -0.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.398'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.27%  '

# Row 41
$ws.Range("E41").Value = '  -1.09%  '

# Row 42
$ws.Range("D42").Value = '3.320.08'
$ws.Range("E42").Value = '  +7.44%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +11.48%  '

# Row 44
$ws.Range("E44").Value = '  +11.20%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0444'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.26%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +16.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.11%  '

# Row 48
$ws.Range("E48").Value = '  +2.43%  '

# Row 49
$ws.Range("E49").Value = '  +8.09%  '

# Row 50
$ws.Range("E50").Value = '  +4.60%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.04%  '
